$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.300.25'
$ws.Range("E2").Value = '  -1.64%  '
$ws.Range("D3").Value = '2.645.77'
$ws.Range("E3").Value = '  +1.02%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '516.65'
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.46'
$ws.Range("E6").Value = '  -1.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.573'
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").Value = '2.677.62'
$ws.Range("E9").Value = '  +2.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.54'
$ws.Range("E10").Value = '  +2.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.105'
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.340'
$ws.Range("E12").Value = '  -0.87%  '
$ws.Range("E13").Value = '  -1.06%  '
$ws.Range("D14").Value = '3.109.60'
$ws.Range("D15").Value = '59.198.89'
$ws.Range("E15").Value = '  -1.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.45'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000139'
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").Value = '2.675.51'
$ws.Range("E18").Value = '  +1.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.62'
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '346.68'
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.54'
$ws.Range("E21").Value = '  +1.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.22'
$ws.Range("E22").Value = '  +0.88%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.03'
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.427'
$ws.Range("E25").Value = '  +1.56%  '
$ws.Range("D26").Value = '2.767.12'
$ws.Range("E26").Value = '  +1.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.993'
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("E28").Value = '  -1.36%  '
$ws.Range("D29").Value = '0.0₃0833'
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.15'
$ws.Range("E30").Value = '  +0.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.67'
$ws.Range("E31").Value = '  +10.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.997'
$ws.Range("E32").Value = '  -0.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.08'
$ws.Range("E33").Value = '  +0.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.58'
$ws.Range("E34").Value = '  -0.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '149.43'
$ws.Range("E35").Value = '  +0.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.04'
$ws.Range("E36").Value = '  +16.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.05'
$ws.Range("E37").Value = '  +1.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.16'
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.873'
$ws.Range("E39").Value = '  -0.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.68'
$ws.Range("E40").Value = '  +0.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.71'
$ws.Range("E41").Value = '  +1.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.44'
$ws.Range("E42").Value = '  -1.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.631'
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '283.47'
$ws.Range("E44").Value = '  -2.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0998'
$ws.Range("E45").Value = '  -0.35%  '
$ws.Range("E46").Value = '  -0.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.77'
$ws.Range("E47").Value = '  +0.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0544'
$ws.Range("E48").Value = '  -1.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.82'
$ws.Range("E49").Value = '  +1.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0233'
$ws.Range("E50").Value = '  -0.83%  '
$ws.Range("E51").Value = '  -1.00%  '
